$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, shifting the former row 107 (guarda lavada) down to 108.
$ws.Rows.Item(107).Insert()

# New row 107 gets the data that row 106 used to have ("1a nueva(o)", 44238).
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 44238
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100114001
$ws.Range("G107").Value = "Papa"
$ws.Range("H107").Value = "Asterix"
$ws.Range("I107").Value = "1a nueva(o)"
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 8000
$ws.Range("L107").Value = 8500
$ws.Range("M107").Value = 8250
$ws.Range("N107").Value = "$/saco 25 kilos"
$ws.Range("O107").Value = "Provincia de Arauco"
$ws.Range("P107").Value = 330
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"

# Row 106 is updated in place: new date and quality label.
$ws.Range("D106").Value = 44448
$ws.Range("I106").Value = "1a (guarda)"
